# Generate Report for Handoff
# The localization file "049022e9-79b8-4129-928f-67c4114b26e1.md" has just been
# handed off again (new Latest Handoff Datetime), so its row moves down (row 3)
# and its Status changes to "Ready for handoff". The file
# "f887603d-a91a-42fb-a406-057eb87fca24.md" (already in sync) moves up to row 2.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Hyperlinks.Delete()

$ws.Range("B2").Value2 = "Handed back: in sync with en-US"
$ws.Range("C2").Value2 = "Handed back: in sync with en-US"

$ws.Range("B3").Value2 = "Ready for handoff"
$ws.Range("C3").Value2 = "Ready for handoff"

$ws.Range("B4").Value2 = "Not to be localized"
$ws.Range("C4").Value2 = "Not to be localized"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/2f7832d980fbb06c9df07fff82cd1115a08570c1/e2e/f887603d-a91a-42fb-a406-057eb87fca24.md", "", "", "f887603d-a91a-42fb-a406-057eb87fca24.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/2f7832d980fbb06c9df07fff82cd1115a08570c1/e2e/049022e9-79b8-4129-928f-67c4114b26e1.md", "", "", "049022e9-79b8-4129-928f-67c4114b26e1.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/2f7832d980fbb06c9df07fff82cd1115a08570c1/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Hyperlinks.Delete()

$ws.Range("B2").Value2 = "Handed back: in sync with en-US"
$ws.Range("D2").Value2 = "2016-03-08 05:10:23"
$ws.Range("G2").Value2 = "2016-03-08 05:11:21"
$ws.Range("H2").Value2 = "Include"

$ws.Range("B3").Value2 = "Ready for handoff"
$ws.Range("D3").Value2 = "2016-03-08 05:12:22"
$ws.Range("G3").Value2 = "2016-03-08 05:11:21"
$ws.Range("H3").Value2 = "Include"

$ws.Range("B4").Value2 = "Not to be localized"
$ws.Range("D4").Value2 = "0001-01-01 00:00:00"
$ws.Range("G4").Value2 = "0001-01-01 00:00:00"
$ws.Range("H4").Value2 = "Ignored"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/2f7832d980fbb06c9df07fff82cd1115a08570c1/e2e/f887603d-a91a-42fb-a406-057eb87fca24.md", "", "", "f887603d-a91a-42fb-a406-057eb87fca24.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f4ad7ddfbd176b6a0f03832959a899a3f7fc5cd3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/f887603d-a91a-42fb-a406-057eb87fca24.35bf874522c2378269b432deb40eec9c5ca62343.zh-cn.xlf", "", "", "f887603d-a91a-42fb-a406-057eb87fca24.35bf874522c2378269b432deb40eec9c5ca62343.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/ad3d806453e57f5208d9e655d463a83eff84ad07/e2e/f887603d-a91a-42fb-a406-057eb87fca24.md", "", "", "f887603d-a91a-42fb-a406-057eb87fca24.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/548eaf1417fe3f7646932a3346874aeb9b382017/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/f887603d-a91a-42fb-a406-057eb87fca24.35bf874522c2378269b432deb40eec9c5ca62343.zh-cn.xlf", "", "", "f887603d-a91a-42fb-a406-057eb87fca24.35bf874522c2378269b432deb40eec9c5ca62343.zh-cn.xlf") | Out-Null

$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/2f7832d980fbb06c9df07fff82cd1115a08570c1/e2e/049022e9-79b8-4129-928f-67c4114b26e1.md", "", "", "049022e9-79b8-4129-928f-67c4114b26e1.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f4ad7ddfbd176b6a0f03832959a899a3f7fc5cd3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/049022e9-79b8-4129-928f-67c4114b26e1.d34f0fc658625c6aa0f66d35039ccd460e31a6fc.zh-cn.xlf", "", "", "049022e9-79b8-4129-928f-67c4114b26e1.d34f0fc658625c6aa0f66d35039ccd460e31a6fc.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/ad3d806453e57f5208d9e655d463a83eff84ad07/e2e/049022e9-79b8-4129-928f-67c4114b26e1.md", "", "", "049022e9-79b8-4129-928f-67c4114b26e1.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/548eaf1417fe3f7646932a3346874aeb9b382017/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/049022e9-79b8-4129-928f-67c4114b26e1.d34f0fc658625c6aa0f66d35039ccd460e31a6fc.zh-cn.xlf", "", "", "049022e9-79b8-4129-928f-67c4114b26e1.d34f0fc658625c6aa0f66d35039ccd460e31a6fc.zh-cn.xlf") | Out-Null

$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/2f7832d980fbb06c9df07fff82cd1115a08570c1/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Hyperlinks.Delete()

$ws.Range("B2").Value2 = "Handed back: in sync with en-US"
$ws.Range("D2").Value2 = "2016-03-08 05:10:36"
$ws.Range("G2").Value2 = "2016-03-08 05:11:38"
$ws.Range("H2").Value2 = "Include"

$ws.Range("B3").Value2 = "Ready for handoff"
$ws.Range("D3").Value2 = "2016-03-08 05:12:32"
$ws.Range("G3").Value2 = "2016-03-08 05:11:38"
$ws.Range("H3").Value2 = "Include"

$ws.Range("B4").Value2 = "Not to be localized"
$ws.Range("D4").Value2 = "0001-01-01 00:00:00"
$ws.Range("G4").Value2 = "0001-01-01 00:00:00"
$ws.Range("H4").Value2 = "Ignored"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/2f7832d980fbb06c9df07fff82cd1115a08570c1/e2e/f887603d-a91a-42fb-a406-057eb87fca24.md", "", "", "f887603d-a91a-42fb-a406-057eb87fca24.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5903ee665c616bbe25c5f56de30e0521eb2d6cab/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/f887603d-a91a-42fb-a406-057eb87fca24.35bf874522c2378269b432deb40eec9c5ca62343.de-de.xlf", "", "", "f887603d-a91a-42fb-a406-057eb87fca24.35bf874522c2378269b432deb40eec9c5ca62343.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/99341443ee52b0e2c8b634b582c00593f9fdc9ed/e2e/f887603d-a91a-42fb-a406-057eb87fca24.md", "", "", "f887603d-a91a-42fb-a406-057eb87fca24.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/4f9f46809c788a37068f9648e53918fb5cb3cd75/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/f887603d-a91a-42fb-a406-057eb87fca24.35bf874522c2378269b432deb40eec9c5ca62343.de-de.xlf", "", "", "f887603d-a91a-42fb-a406-057eb87fca24.35bf874522c2378269b432deb40eec9c5ca62343.de-de.xlf") | Out-Null

$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/2f7832d980fbb06c9df07fff82cd1115a08570c1/e2e/049022e9-79b8-4129-928f-67c4114b26e1.md", "", "", "049022e9-79b8-4129-928f-67c4114b26e1.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5903ee665c616bbe25c5f56de30e0521eb2d6cab/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/049022e9-79b8-4129-928f-67c4114b26e1.d34f0fc658625c6aa0f66d35039ccd460e31a6fc.de-de.xlf", "", "", "049022e9-79b8-4129-928f-67c4114b26e1.d34f0fc658625c6aa0f66d35039ccd460e31a6fc.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/99341443ee52b0e2c8b634b582c00593f9fdc9ed/e2e/049022e9-79b8-4129-928f-67c4114b26e1.md", "", "", "049022e9-79b8-4129-928f-67c4114b26e1.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/4f9f46809c788a37068f9648e53918fb5cb3cd75/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/049022e9-79b8-4129-928f-67c4114b26e1.d34f0fc658625c6aa0f66d35039ccd460e31a6fc.de-de.xlf", "", "", "049022e9-79b8-4129-928f-67c4114b26e1.d34f0fc658625c6aa0f66d35039ccd460e31a6fc.de-de.xlf") | Out-Null

$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/2f7832d980fbb06c9df07fff82cd1115a08570c1/.localization-config", "", "", ".localization-config") | Out-Null

Write-Host "Done."
